# 439-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-B-EarlyRePayment-Loanproduct4
# "code refactoring and loan accounting and charges added"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoan_Input
$ws2 = $wb.Worksheets.Item(2)   # ProductLoan_Output

# --- Rename the product (shared by both sheets' B1 cell) ---
$ws1.Range("B1").Value = "439-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-B-EarlyRePayment"
$ws2.Range("B1").Value = "439-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-B-EarlyRePayment"

# --- shortname switches from text "kar6" to the numeric code 439 ---
$ws1.Range("B3").Value = 439

# --- maximumtranchecount goes from 12 down to 1 ---
$ws1.Range("B11").Value = 1

# --- Append the new loan-accounting / charges rows (29-40) ---
# First lay down the style for the new rows by copying the existing
# label/value style pairing (row 10 has the same A/B style combo we need).
$ws1.Range("A10:B10").Copy($ws1.Range("A29:B40"))

# Column A labels (field names) written first, top-to-bottom ...
$ws1.Range("A29").Value = "fundsource"
$ws1.Range("A30").Value = "loanprotfolio"
$ws1.Range("A31").Value = "interestreceivable"
$ws1.Range("A32").Value = "penaltiesreceivable"
$ws1.Range("A33").Value = "transferinsuspense"
$ws1.Range("A34").Value = "feesreceivable"
$ws1.Range("A35").Value = "incomefrominterest"
$ws1.Range("A36").Value = "incomefrompenalties"
$ws1.Range("A37").Value = "incomefromfees"
$ws1.Range("A38").Value = "incomefromrecoveryrepayments"
$ws1.Range("A39").Value = "loseswrittenoff"
$ws1.Range("A40").Value = "overpaymentliability"

# ... then column B values (accounting labels), top-to-bottom
$ws1.Range("B29").Value = "Cash"
$ws1.Range("B30").Value = "Loan portfolio "
$ws1.Range("B31").Value = "Interest Receivable "
$ws1.Range("B32").Value = "Penalties Receivable "
$ws1.Range("B33").Value = "Transfer in Suspence "
$ws1.Range("B34").Value = "Fees Receivable"
$ws1.Range("B35").Value = "Income from interest"
$ws1.Range("B36").Value = "Income from penalties"
$ws1.Range("B37").Value = "Income from fees"
$ws1.Range("B38").Value = "Income from recovery repayments"
$ws1.Range("B39").Value = "Losses Writtenoff "
$ws1.Range("B40").Value = "Overpayment Liability"

# --- View-state tweaks: activate the Output sheet, scroll the Input sheet ---
$ws1.Range("A26").Select()
$ws2.Range("E15").Select()
$ws2.Activate()
